$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.064.85'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +4.65%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.236.91'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +2.94%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.93'
$ws.Range("D5").Style = "Normal"

$ws.Range("E6").Value = '  +0.82%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '75.68'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  +7.88%  '

$ws.Range("E8").Value = '  -0.12%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.617'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  +6.65%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.11'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  +2.41%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0929'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  +0.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.47'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  +0.29%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.97'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  +3.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.102'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +0.85%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.568.70'
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.71'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +5.71%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.230.32'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  +3.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.812'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  +0.38%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.939.16'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  +4.90%  '

$ws.Range("E20").Value = '  +3.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.13'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  +1.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.01'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  +1.30%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.41'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  +3.57%  '

$ws.Range("B24").Value = 'BitcoinCash'

$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '230.61'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  +2.25%  '

$ws.Range("B25").Value = 'ImmutableX'

$ws.Range("C25").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.20'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  +13.02%  '

$ws.Range("E26").Value = '  -0.10%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.99'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  +0.95%  '

$ws.Range("E28").Value = '  -4.50%  '

$ws.Range("E29").Value = '  +2.30%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.04'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  +4.36%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '37.38'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  +20.88%  '

$ws.Range("E32").Value = '  -3.45%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.35'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  +2.16%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0794'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  +3.17%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.37'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  +3.86%  '

$ws.Range("E36").Value = '  +1.21%  '

$ws.Range("E37").Value = '  +7.46%  '

$ws.Range("E38").Value = '  +5.43%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0333'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +16.69%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '13.13'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  +5.72%  '

$ws.Range("E41").Value = '  +3.16%  '

$ws.Range("E42").Value = '  +3.20%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.200'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  +5.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '60.09'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  -0.28%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '105.01'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  +6.78%  '

$ws.Range("E46").Value = '  +3.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0994'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  +2.40%  '

$ws.Range("E48").Value = '  +4.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.444'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +19.10%  '

$ws.Range("E50").Value = '  +1.42%  '

$ws.Range("E51").Value = '  +1.97%  '
